# Add two new columns (I = "I0", J = "IF") with numeric data for rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns, matching the style of the existing headers (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2 through 39
$iValues = @(9,3,7,5,11,7,9,8,7,7,8,7,5,6,6,6,7,4,7,7,8,7,7,8,6,8,6,5,3,9,8,2,8,8,6,9,6,6)
$jValues = @(9,3,7,5,11,7,9,8,7,7,8,8,5,7,6,6,8,5,8,7,8,7,7,8,6,8,6,5,3,9,8,2,8,8,6,9,6,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
